# Corrección de links a código .qmd desde calendario.
# Replace the ".Rmd" code-file links/paths with ".qmd" across the calendar
# sheets (Sheet1 col G, Sheet2 col B), then restore the two sheets' active
# selections as left by the author after the edit.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($r + $rowOffset, $c + $colOffset)
            $val = $cell.Value2
            if ($val -ne $null -and $val -is [string] -and $val.Contains(".Rmd")) {
                $cell.Value2 = $val.Replace(".Rmd", ".qmd")
            }
        }
    }
}

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("E20").Select()

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Activate()
$sheet2.Range("H15").Select()

$sheet1.Activate()
